# Commit: "time between injuries complete"
# Expand the abbreviated level value "JV" to "Junior Varsity" throughout the
# data sheet (columns level_1 / level_2), matching only whole-cell contents
# so no other text is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlWhole = 1 ensures only cells whose entire content is exactly "JV" are
# replaced (not substrings inside other values).
$ws.Cells.Replace("JV", "Junior Varsity", 1)
